$wb = $excel.ActiveWorkbook

# --- 1. Update status text "Ready for handoff" -> "In Translation" ---
# Overview sheet: columns E (zh-cn) and F (de-de), rows 2-3
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# zh-cn sheet: column C (Status), rows 2-3
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

# de-de sheet: column C (Status), rows 2-3
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- 2. Narrow the "Status" columns ---
# Overview: columns E & F (zh-cn / de-de status columns)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn: column C (Status)
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# de-de: column C (Status)
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
